# Fix Training Data Issue: the "Date" column (BF) was stored as the sheet's
# own file-name-derived label ("4-20-2013-14") instead of an ISO date
# ("2014-04-20"). NBA stats for games played late in the evening were
# attributed to the wrong calendar day, so the literal text in every data
# row needs to change from "4-20-2013-14" to "2014-04-20".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$oldText = "4-20-2013-14"
$newText = "2014-04-20"

$ur = $ws.UsedRange
$lastRow = $ur.Rows.Count
$lastCol = $ur.Columns.Count

# Locate the "Date" header column dynamically (it's BF == column 58 in this
# workbook, but search for it so the script isn't brittle to layout drift).
$dateCol = 0
for ($c = 1; $c -le $lastCol; $c++) {
    if ($ws.Cells.Item(1, $c).Text -eq "Date") {
        $dateCol = $c
    }
}
if ($dateCol -eq 0) { $dateCol = 58 }

# Walk every data row and replace the old literal date-ish text with the
# corrected ISO date. Cells are forced to Text format first so Excel does
# not auto-convert the "2014-04-20" literal into a date serial number —
# the column must keep holding plain text, exactly as it did before.
for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, $dateCol)
    if ($cell.Text -eq $oldText) {
        $cell.NumberFormat = "@"
        $cell.Value = $newText
        $cell.ClearFormats()
    }
}
